# Replace the single "m:'prefix\nsuffix'" field (fldChar begin/instrText.../fldChar end)
# with literal text runs "{", "m", ":'", "prefix", "\n", "suffix", "'", "}" while keeping
# the _GoBack bookmark in its original position (between "prefix" and "\n").
$d = $word.ActiveDocument

# Locate the paragraph that hosts the field (robust to paragraph index/content drift).
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $targetPara = $p
        break
    }
}

if ($null -eq $targetPara) {
    throw "Could not find the paragraph containing the field to rewrite."
}

$newParaXml = '<w:p>' +
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m</w:t></w:r>' +
    '<w:r><w:t>:''</w:t></w:r>' +
    '<w:r><w:t>prefix</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>\n</w:t></w:r>' +
    '<w:r><w:t>suffix</w:t></w:r>' +
    '<w:r><w:t>''</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
    '</w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $newParaXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$targetPara.Range.InsertXML($xml) | Out-Null

Write-Output "Field rewritten to literal token text."
